$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the user-story text (column C) — fix wording / typos and
#    shorten a couple of the sentences.
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "utilisateur connecté"
$ws.Range("C4").Value = "ajouter un contact"
$ws.Range("C7").Value = "L'utilisateur ajoute un contact afin de pouvoir lui envoyer des messages"
$ws.Range("C9").Value = 'clique sur le boutton "+"'
$ws.Range("C10").Value = "je peux ajouter un contact"

# ---------------------------------------------------------------------
# 2. Give row 7 ("Scénario") its own thick-bordered box, matching the
#    treatment row 5 ("Afin de") already has, by adding a medium bottom
#    border under B7:C7 (it already has thick left/right + medium top).
# ---------------------------------------------------------------------
$ws.Range("B7").Borders.Item(9).Weight = -4138
$ws.Range("C7").Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------
# 3. Make row 7 look like the other "closed" rows: bottom border + taller.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 24
}

# ---------------------------------------------------------------------
# 4. View tweaks: zoom in a bit and move the selection off the table.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("C17").Select()
